$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 6: tcid=5, priority=P2, username=mmanubrolu, password=mmanubrolu
$ws.Cells.Item(6, 1).Value = "'5"
$ws.Cells.Item(6, 2).Value = "P2"
$ws.Cells.Item(6, 3).Value = "mmanubrolu"
$ws.Cells.Item(6, 4).Value = "mmanubrolu"

$ws.Range("H11").Select()
